# ERP ID logic implementation
# Renames the FF_0301xx "ERP ID" test values to FF_0303xx across the
# three worksheets that reference them, and updates the active
# sheet/selection bookmarks left behind by the editor.

$wb = $excel.ActiveWorkbook

$wsLogin   = $wb.Worksheets.Item("loginTest")
$wsAccount = $wb.Worksheets.Item("accountCreatonTest")
$wsProfile = $wb.Worksheets.Item("profileUpdateTest")

# --- accountCreatonTest -------------------------------------------------
$wsAccount.Range("E2").Value  = "FF_030316aa"
$wsAccount.Range("H2").Value  = "FF_030316aa"
$wsAccount.Range("I2").Value  = "FF_030316aa"
$wsAccount.Range("J2").Value  = "FF_030316aa"

$wsAccount.Range("E3").Value  = "FF_030316ab"
$wsAccount.Range("H3").Value  = "FF_030316ab"
$wsAccount.Range("I3").Value  = "FF_030316ab"
$wsAccount.Range("J3").Value  = "FF_030316ab"

$wsAccount.Range("E4").Value  = "FF_030316ac"
$wsAccount.Range("H4").Value  = "FF_030316ac"
$wsAccount.Range("I4").Value  = "FF_030316ac"
$wsAccount.Range("J4").Value  = "FF_030316ac"

$wsAccount.Range("E5").Value  = "FF_030316ad"
$wsAccount.Range("H5").Value  = "FF_030316ad"
$wsAccount.Range("I5").Value  = "FF_030316ad"
$wsAccount.Range("J5").Value  = "FF_030316ad"

$wsAccount.Range("E6").Value  = "FF_030316ae"
$wsAccount.Range("H6").Value  = "FF_030316ae"
$wsAccount.Range("I6").Value  = "FF_030316ae"
$wsAccount.Range("J6").Value  = "FF_030316ae"

$wsAccount.Range("E7").Value  = "FF_030316af"
$wsAccount.Range("H7").Value  = "FF_030316af"
$wsAccount.Range("I7").Value  = "FF_030316af"
$wsAccount.Range("J7").Value  = "FF_030316af"

$wsAccount.Range("E8").Value  = "FF_030316ag"
$wsAccount.Range("H8").Value  = "FF_030316ag"
$wsAccount.Range("I8").Value  = "FF_030316ag"
$wsAccount.Range("J8").Value  = "FF_030316ag"

$wsAccount.Range("E9").Value  = "FF_030316ah"
$wsAccount.Range("H9").Value  = "FF_030316ah"
$wsAccount.Range("I9").Value  = "FF_030316ah"
$wsAccount.Range("J9").Value  = "FF_030316ah"

$wsAccount.Range("E10").Value = "FF_030316ai"
$wsAccount.Range("H10").Value = "FF_030316ai"
$wsAccount.Range("I10").Value = "FF_030316ai"
$wsAccount.Range("J10").Value = "FF_030316ai"

$wsAccount.Range("E11").Value = "FF_030316aj"
$wsAccount.Range("H11").Value = "FF_030316aj"
$wsAccount.Range("I11").Value = "FF_030316aj"
$wsAccount.Range("J11").Value = "FF_030316aj"

# --- loginTest ---------------------------------------------------------
$wsLogin.Range("B2").Value = "FF_030316ac"
$wsLogin.Range("C2").Value = "FF_030316ac"

# --- profileUpdateTest ---------------------------------------------------
$wsProfile.Range("E2").Value = "FF_030316aa"

# --- Selections / active sheet ------------------------------------------
# Leave the same bookmark behaviour the author ended up with: profile and
# account sheets keep a plain selection, while loginTest becomes the
# active (tabSelected) sheet with the final selection.
$wsProfile.Range("A2").Select()
$wsAccount.Range("E2").Select()
$wsLogin.Range("B6").Select()
